# Apply "Different code for SiPM and PMT/GEM" edits to CoG_test.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update the input data cells (D4, E4, E5, D6). Dependent formulas
# (J4, K4, R2) recalculate automatically.
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 20
$ws.Range("D6").Value = 4

# Force a full recalculation so cached formula values stay in sync.
$excel.CalculateFullRebuild()

# Move the active selection from G9 to S2, matching the saved view state.
$ws.Activate()
$ws.Range("S2").Select()
